# Noted the completed GDML loading tests
# ---------------------------------------------------------------------------
# Applies the Test.xlsx edit: extends the "Solids" sheet with the full list
# of GDML solid types (and which stages of the pipeline support them),
# tidies up the "Other" sheet (removing the now-redundant "multiunion" row
# and re-selecting "Solids" as the active tab), and widens the first
# column on "Physical" to fit its longer header text.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$xlCenter = -4108

# ---------------------------------------------------------------------------
# "Other" sheet: touch it first (selecting a range here would otherwise make
# it the active tab again after we activate "Solids" below).
# ---------------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("Other")
$wsOther.Columns.Item(1).ColumnWidth = 20.42
$wsOther.Range("B2").Select()

# ---------------------------------------------------------------------------
# "Physical" sheet: only a first-column width tweak (bestFit for new header).
# ---------------------------------------------------------------------------
$wsPhys = $wb.Worksheets.Item("Physical")
$wsPhys.Columns.Item(1).ColumnWidth = 14.25

# ---------------------------------------------------------------------------
# "Solids" sheet: the main content of the edit.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Solids")
$ws.Columns.Item(1).ColumnWidth = 18.09

$names = @("Box","Tube","Cut Tube","Cone","Para","Trd","Trap","Sphere","Orb","Torus","Polycone","Generic polycone","Polyhedra","Generic polyhedra","Ellipical Tube","Ellipsoid","Elliptical cone","Paraboloid","Hype","Tet","Extrusion solid","Twisted box","Twisted trap","Twisted trd","Arbitrary trap","Tessellated solid","Union","Subtraction","Intersection","Multiunion")
$hasY  = @($true,$true,$true,$true,$true,$false,$false,$true,$true,$true,$true,$false,$true,$false,$true,$true,$true,$true,$true,$true,$false,$false,$false,$false,$false,$false,$true,$true,$true,$false)

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $names[$i]
    if ($hasY[$i]) {
        $ws.Cells.Item($r, 2).Value = "Y"
    }
    $ws.Cells.Item($r, 2).HorizontalAlignment = $xlCenter
}

# Totals row (per-column count of ticked solids).
$ws.Range("A32").Value = "Total"
$ws.Range("B32").Formula = "=COUNTA(B2:B31)"
$ws.Range("C32:D32").Formula = "=COUNTA(C2:C31)"
$ws.Range("B32:D32").HorizontalAlignment = $xlCenter

# Grand total row.
$ws.Range("A35").Value = "Total solids"
$ws.Range("A35").Font.Bold = $true
$ws.Range("B35").Formula = "=COUNTBLANK(B1:B31) + COUNTA(B2:B31)"
$ws.Range("B35").Font.Bold = $true
$ws.Range("B35").HorizontalAlignment = $xlCenter

$ws.Range("A36").Font.Bold = $true
$ws.Range("B36").Font.Bold = $true
$ws.Range("B36").HorizontalAlignment = $xlCenter

# Re-select "Solids" as the active tab, as in the authored workbook.
$ws.Select()
$ws.Range("B14").Select()
